# Auto-generated edit script applying the Sargatanas_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4500
$ws.Range("J64").Value = 4500
$ws.Range("L64").Value = 4500
$ws.Range("N64").Value = -4996

$ws.Range("H67").Value = 4500
$ws.Range("J67").Value = 4500
$ws.Range("L67").Value = 4500
$ws.Range("N67").Value = -6216

$ws.Range("H86").Value = 97225544
$ws.Range("I86").Value = 125002760
$ws.Range("K86").Value = 125002760
$ws.Range("M86").Value = -125001637

$ws.Range("H87").Value = 51998.4
$ws.Range("J87").Value = 51998.4
$ws.Range("L87").Value = 51998.4
$ws.Range("N87").Value = -54494.4

$ws.Range("H89").Value = 97225544
$ws.Range("I89").Value = 125002760
$ws.Range("K89").Value = 625013800
$ws.Range("M89").Value = -625008184

$ws.Range("H90").Value = 51998.4
$ws.Range("J90").Value = 51998.4
$ws.Range("L90").Value = 155995.2
$ws.Range("N90").Value = -168475.2

$ws.Range("H100").Value = 2367.9473
$ws.Range("I100").Value = 1776.3077
$ws.Range("J100").Value = 3649.8333
$ws.Range("K100").Value = 1776.3077
$ws.Range("L100").Value = 3649.8333
$ws.Range("M100").Value = -1235.3077
$ws.Range("N100").Value = -4731.8333

$ws.Range("H135").Value = 1429460.9
$ws.Range("I135").Value = 1429460.9
$ws.Range("K135").Value = 12865148.1
$ws.Range("M135").Value = -12862613.1

$ws.Range("H137").Value = 3384.182
$ws.Range("I137").Value = 4207
$ws.Range("J137").Value = 3000.2
$ws.Range("K137").Value = 12621
$ws.Range("L137").Value = 9000.599999999999
$ws.Range("M137").Value = -10071
$ws.Range("N137").Value = -14100.6

$ws.Range("H138").Value = 6823.2173
$ws.Range("I138").Value = 1308.75
$ws.Range("J138").Value = 19427.715
$ws.Range("K138").Value = 3926.25
$ws.Range("L138").Value = 58283.145
$ws.Range("M138").Value = 1213.75
$ws.Range("N138").Value = -68563.145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2852.6956
$ws.Range("I2").Value = 1333.7858
$ws.Range("K2").Value = 1333.7858
$ws.Range("M2").Value = -1220.7858

$ws.Range("H30").Value = 1983.3334
$ws.Range("I30").Value = 1550
$ws.Range("J30").Value = 2850
$ws.Range("K30").Value = 1550
$ws.Range("L30").Value = 2850
$ws.Range("M30").Value = -1400
$ws.Range("N30").Value = -3150

$ws.Range("H57").Value = 4745.4546
$ws.Range("I57").Value = 4745.4546
$ws.Range("K57").Value = 4745.4546
$ws.Range("M57").Value = -4261.4546

$ws.Range("H116").Value = 2852.6956
$ws.Range("I116").Value = 1333.7858
$ws.Range("K116").Value = 1333.7858
$ws.Range("M116").Value = 960.2141999999999

$ws.Range("H122").Value = 2333.1072
$ws.Range("I122").Value = 1117.6364
$ws.Range("K122").Value = 3352.9092
$ws.Range("M122").Value = -902.9092000000001

$ws.Range("H126").Value = 5241.75
$ws.Range("I126").Value = 5241.75
$ws.Range("K126").Value = 15725.25
$ws.Range("M126").Value = -13255.25

$ws.Range("H132").Value = 2709
$ws.Range("I132").Value = 1516.5
$ws.Range("K132").Value = 4549.5
$ws.Range("M132").Value = -2019.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2852.6956
$ws.Range("I3").Value = 1333.7858
$ws.Range("K3").Value = 1333.7858
$ws.Range("M3").Value = -1219.7858

$ws.Range("H20").Value = 11113474
$ws.Range("I20").Value = 15153384
$ws.Range("K20").Value = 15153384
$ws.Range("M20").Value = -15153137

$ws.Range("H107").Value = 40182270
$ws.Range("I107").Value = 45002744
$ws.Range("K107").Value = 45002744
$ws.Range("M107").Value = -45000824

$ws.Range("H125").Value = 51548.5
$ws.Range("J125").Value = 51548.5
$ws.Range("L125").Value = 51548.5
$ws.Range("N125").Value = -61388.5

$ws.Range("H128").Value = 3839.5
$ws.Range("I128").Value = 3839.5
$ws.Range("K128").Value = 11518.5
$ws.Range("M128").Value = -9028.5

$ws.Range("H134").Value = 4504.65
$ws.Range("I134").Value = 1693.027
$ws.Range("K134").Value = 5079.081
$ws.Range("M134").Value = -2544.081

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9366.143
$ws.Range("I31").Value = 4510.875
$ws.Range("J31").Value = 10804.741
$ws.Range("K31").Value = 4510.875
$ws.Range("L31").Value = 10804.741
$ws.Range("M31").Value = -4215.875
$ws.Range("N31").Value = -11394.741

$ws.Range("H34").Value = 9366.143
$ws.Range("I34").Value = 4510.875
$ws.Range("J34").Value = 10804.741
$ws.Range("K34").Value = 4510.875
$ws.Range("L34").Value = 10804.741
$ws.Range("M34").Value = -4308.875
$ws.Range("N34").Value = -11208.741

$ws.Range("H76").Value = 4953.857
$ws.Range("I76").Value = 4953.857
$ws.Range("K76").Value = 4953.857
$ws.Range("M76").Value = -4638.857

$ws.Range("H79").Value = 4953.857
$ws.Range("I79").Value = 4953.857
$ws.Range("K79").Value = 4953.857
$ws.Range("M79").Value = -3861.857

$ws.Range("H132").Value = 5271
$ws.Range("I132").Value = 1906
$ws.Range("K132").Value = 5718
$ws.Range("M132").Value = -3188

$ws.Range("H134").Value = 4670.5107
$ws.Range("I134").Value = 1279.871
$ws.Range("K134").Value = 3839.613
$ws.Range("M134").Value = -1304.613

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 60501336
$ws.Range("I4").Value = 136111500
$ws.Range("J4").Value = 13201.1
$ws.Range("K4").Value = 408334500
$ws.Range("L4").Value = 39603.3
$ws.Range("M4").Value = -408334388
$ws.Range("N4").Value = -39827.3

$ws.Range("H34").Value = 4162.2104
$ws.Range("I34").Value = 519.8
$ws.Range("J34").Value = 5463.0713
$ws.Range("K34").Value = 1559.4
$ws.Range("L34").Value = 16389.2139
$ws.Range("M34").Value = -1475.4
$ws.Range("N34").Value = -16557.2139

$ws.Range("H114").Value = 592.3
$ws.Range("I114").Value = 203.28572
$ws.Range("K114").Value = 609.85716
$ws.Range("M114").Value = 2644.14284

$ws.Range("H138").Value = 57272.316
$ws.Range("J138").Value = 9009.333000000001
$ws.Range("L138").Value = 27027.999
$ws.Range("N138").Value = -37307.999

$ws.Range("H141").Value = 9677.416999999999
$ws.Range("I141").Value = 3732.7144
$ws.Range("K141").Value = 11198.1432
$ws.Range("M141").Value = -6018.143199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8596.958000000001
$ws.Range("I70").Value = 7739.8
$ws.Range("J70").Value = 9209.214
$ws.Range("K70").Value = 7739.8
$ws.Range("L70").Value = 9209.214
$ws.Range("M70").Value = -7469.8
$ws.Range("N70").Value = -9749.214

$ws.Range("H73").Value = 8596.958000000001
$ws.Range("I73").Value = 7739.8
$ws.Range("J73").Value = 9209.214
$ws.Range("K73").Value = 7739.8
$ws.Range("L73").Value = 9209.214
$ws.Range("M73").Value = -6803.8
$ws.Range("N73").Value = -11081.214

$ws.Range("H122").Value = 2749809.5
$ws.Range("I122").Value = 3573415.8
$ws.Range("K122").Value = 10720247.4
$ws.Range("M122").Value = -10717797.4

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H132").Value = 4132.1934
$ws.Range("I132").Value = 1669.2174
$ws.Range("J132").Value = 11213.25
$ws.Range("K132").Value = 5007.6522
$ws.Range("L132").Value = 33639.75
$ws.Range("M132").Value = -2477.6522
$ws.Range("N132").Value = -38699.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1652.6428
$ws.Range("J46").Value = 2529.125
$ws.Range("L46").Value = 2529.125
$ws.Range("N46").Value = -2905.125

$ws.Range("H55").Value = 32258398
$ws.Range("I55").Value = 62500052
$ws.Range("J55").Value = 632.93335
$ws.Range("K55").Value = 62500052
$ws.Range("L55").Value = 632.93335
$ws.Range("M55").Value = -62499879
$ws.Range("N55").Value = -978.93335

$ws.Range("H125").Value = 51549
$ws.Range("J125").Value = 51549
$ws.Range("L125").Value = 51549
$ws.Range("N125").Value = -61389

$ws.Range("H132").Value = 10876083
$ws.Range("I132").Value = 27780238
$ws.Range("K132").Value = 83340714
$ws.Range("M132").Value = -83338184

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5702.5806
$ws.Range("I122").Value = 5237.143
$ws.Range("K122").Value = 15711.429
$ws.Range("M122").Value = -13261.429

$ws.Range("H126").Value = 2673.1
$ws.Range("I126").Value = 1680.7333
$ws.Range("K126").Value = 5042.199900000001
$ws.Range("M126").Value = -2572.199900000001

$ws.Range("H132").Value = 29435386
$ws.Range("I132").Value = 62514252
$ws.Range("J132").Value = 31948.777
$ws.Range("K132").Value = 187542756
$ws.Range("L132").Value = 95846.33099999999
$ws.Range("M132").Value = -187540226
$ws.Range("N132").Value = -100906.331
